$d = $word.ActiveDocument

# Locate the run of text "Мельников " (with its trailing space) right
# after which the supervisor's initials need to be appended.
$r = $d.Content
$found = $r.Find.Execute("Мельников ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Collapse to the end of the match and select it so formatting applied
# through the Selection object (needed to get both the ascii/hAnsi/cs
# font names and both sz/szCs sizes serialized) lands on the new text.
$r.Collapse(0)
$r.Select()
$r.InsertAfter("Д.А.")

# Give the newly inserted run the same formatting used throughout this
# paragraph/document: Times New Roman, 14pt (28 half-points).
$sel = $word.Selection
$sel.Font.Name = "Times New Roman"
$sel.Font.NameBi = "Times New Roman"
$sel.Font.Size = 14
$sel.Font.SizeBi = 14
